# Updated cryptos list prices (Price column D) and Volume(1h) percentages (column E).
# Price values keep their original text formatting (e.g. "1.001", "28.145.30"),
# so they are written with a leading apostrophe to force Excel to store them as text
# instead of auto-converting them to numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = '28.145.30'; Volume = '  -3.19%  ' }
    @{ Row = 3; Price = '1.926.18'; Volume = '  -2.35%  ' }
    @{ Row = 4; Price = '1.001'; Volume = '  -1.03%  ' }
    @{ Row = 5; Price = '330.43'; Volume = '  +0.38%  ' }
    @{ Row = 6; Price = '1.001'; Volume = '  -0.95%  ' }
    @{ Row = 7; Price = '0.4726'; Volume = '  -4.79%  ' }
    @{ Row = 8; Price = '0.4056'; Volume = '  -3.63%  ' }
    @{ Row = 9; Price = '53.13'; Volume = '  -0.68%  ' }
    @{ Row = 10; Price = '0.08428'; Volume = '  -9.03%  ' }
    @{ Row = 11; Price = '1.048'; Volume = '  -4.65%  ' }
    @{ Row = 12; Price = '22.29'; Volume = '  -2.43%  ' }
    @{ Row = 13; Price = '1.909.95'; Volume = '  -3.29%  ' }
    @{ Row = 14; Price = '7.509'; Volume = '  -4.91%  ' }
    @{ Row = 15; Price = '6.101'; Volume = '  -5.43%  ' }
    @{ Row = 16; Price = '1.003'; Volume = '  -0.98%  ' }
    @{ Row = 17; Price = '90.53'; Volume = '  -1.38%  ' }
    @{ Row = 18; Price = '0.00001067'; Volume = '  -3.68%  ' }
    @{ Row = 19; Price = '0.06596'; Volume = '  -1.83%  ' }
    @{ Row = 20; Price = '18.10'; Volume = '  -5.45%  ' }
    @{ Row = 21; Price = '1.002'; Volume = '  -0.65%  ' }
    @{ Row = 22; Price = '5.751'; Volume = '  -3.46%  ' }
    @{ Row = 23; Price = '28.156.03'; Volume = '  -3.21%  ' }
    @{ Row = 24; Price = '11.38'; Volume = '  -4.81%  ' }
    @{ Row = 25; Price = '2.285'; Volume = '  +0.88%  ' }
    @{ Row = 26; Price = '2.160.02'; Volume = '  -2.11%  ' }
    @{ Row = 27; Price = '154.19'; Volume = '  -0.90%  ' }
    @{ Row = 28; Price = '20.09'; Volume = '  -3.00%  ' }
    @{ Row = 29; Price = '2.152'; Volume = '  -4.81%  ' }
    @{ Row = 30; Price = '5.744'; Volume = '  -8.62%  ' }
    @{ Row = 31; Price = '123.75'; Volume = '  -2.65%  ' }
    @{ Row = 32; Price = '0.9785'; Volume = '  -6.47%  ' }
    @{ Row = 33; Price = '0.09616'; Volume = '  -2.36%  ' }
    @{ Row = 34; Price = '1.455'; Volume = '  -3.49%  ' }
    @{ Row = 35; Price = '5.569'; Volume = '  -4.28%  ' }
    @{ Row = 36; Price = '3.633'; Volume = '  -2.79%  ' }
    @{ Row = 37; Price = '9.017'; Volume = '  -0.17%  ' }
    @{ Row = 38; Price = '0.02316'; Volume = '  -4.33%  ' }
    @{ Row = 39; Price = '0.06178'; Volume = '  -3.82%  ' }
    @{ Row = 40; Price = '1.235'; Volume = '  -6.88%  ' }
    @{ Row = 41; Price = '0.6168'; Volume = '  -4.67%  ' }
    @{ Row = 42; Price = '11.08'; Volume = '  -3.66%  ' }
    @{ Row = 43; Price = '1.001'; Volume = '  -0.81%  ' }
    @{ Row = 44; Price = '0.1906'; Volume = '  -4.79%  ' }
    @{ Row = 45; Price = '1.306'; Volume = '  -4.52%  ' }
    @{ Row = 46; Price = '0.5886'; Volume = '  -5.08%  ' }
    @{ Row = 47; Price = '12.80'; Volume = '  -3.90%  ' }
    @{ Row = 48; Price = '2.036'; Volume = '  -6.69%  ' }
    @{ Row = 49; Price = '3.474'; Volume = '  -0.31%  ' }
    @{ Row = 50; Price = '0.06844'; Volume = '  -1.97%  ' }
    @{ Row = 51; Price = '110.27'; Volume = '  -2.62%  ' }
)

foreach ($u in $updates) {
    $priceCell = "D" + $u.Row
    $volumeCell = "E" + $u.Row
    # Leading apostrophe forces Excel to store the price as literal text
    # (matches the source data, which keeps thousands separators as dots,
    # e.g. "28.145.30", and would otherwise be auto-converted to a number).
    $ws.Range($priceCell).Value = "'" + $u.Price
    $ws.Range($volumeCell).Value = $u.Volume
}
